$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 51, shifting existing rows 51-56 down to 52-57
$ws.Rows.Item(51).Insert()

# Populate the new row 51 with the new data record
$ws.Cells.Item(51, 1).Value = 11
$ws.Cells.Item(51, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(51, 3).Value = "Bíobío"
$ws.Cells.Item(51, 4).Value = 44516
$ws.Cells.Item(51, 4).Style = $ws.Cells.Item(52, 4).Style
$ws.Cells.Item(51, 5).Value = 8
$ws.Cells.Item(51, 6).Value = 100112021
$ws.Cells.Item(51, 7).Value = "Ají"
$ws.Cells.Item(51, 8).Value = "Americana (o)"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 550
$ws.Cells.Item(51, 11).Value = 14000
$ws.Cells.Item(51, 12).Value = 15000
$ws.Cells.Item(51, 13).Value = 14545
$ws.Cells.Item(51, 14).Value = "`$/caja 12 kilos"
$ws.Cells.Item(51, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(51, 16).Value = 1212
$ws.Cells.Item(51, 17).Value = 12
$ws.Cells.Item(51, 18).Value = "Hortaliza"
